$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for B2:I21 (reshuffled with new random seed)
$data = @(
    @(0.0,0.02,0.04,-0.05,-0.0,0.03,0.06,-0.07),
    @(0.01,0.04,0.08,-0.09,-0.0,0.06,0.13,-0.16),
    @(-0.0,0.02,0.03,-0.04,-0.0,0.03,0.06,-0.08),
    @(-0.0,0.01,0.03,-0.04,-0.0,0.02,0.05,-0.05),
    @(0.0,0.03,0.07,-0.07,-0.0,0.04,0.09,-0.12),
    @(0.0,0.02,0.05,-0.05,-0.0,0.03,0.07,-0.08),
    @(-0.0,0.02,0.05,-0.06,-0.0,0.03,0.09,-0.09),
    @(-0.0,0.03,0.06,-0.08,-0.01,0.04,0.09,-0.12),
    @(-0.0,0.01,0.03,-0.04,-0.0,0.02,0.05,-0.06),
    @(-0.0,0.02,0.06,-0.05,-0.0,0.03,0.07,-0.09),
    @(-0.0,0.02,0.04,-0.04,-0.0,0.03,0.05,-0.11),
    @(-0.0,0.02,0.04,-0.05,-0.01,0.03,0.07,-0.1),
    @(0.0,0.01,0.04,-0.04,-0.0,0.02,0.05,-0.06),
    @(0.0,0.02,0.04,-0.04,-0.0,0.03,0.06,-0.06),
    @(0.0,0.02,0.03,-0.04,-0.0,0.02,0.07,-0.05),
    @(-0.0,0.01,0.03,-0.03,-0.0,0.02,0.05,-0.07),
    @(-0.0,0.03,0.08,-0.06,-0.01,0.04,0.09,-0.13),
    @(-0.0,0.02,0.05,-0.07,-0.01,0.04,0.07,-0.1),
    @(0.0,0.04,0.09,-0.08,-0.01,0.05,0.11,-0.13),
    @(0.0,0.03,0.08,-0.08,-0.0,0.05,0.13,-0.12)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $c = $j + 2
        $ws.Cells.Item($r, $c).Value = $rowVals[$j]
    }
}